$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "string1"
$ws.Range("F9").Value = 456
$ws.Range("H9").Value = 106704

$ws.Range("H31").Value = 106704
$ws.Range("H32").Value = 105636.96
$ws.Range("H33").Value = 118313.3952
